# Reorders rows 2-19 (Fecha + the dependent data columns D..T) to match the
# "weekly" resequencing described in the commit. Columns A-C (Mercado ID,
# Mercado, Region) and E-L (Codreg..Calidad) are constant across every row
# in this subset, so only D and M..T need to move.
#
# Mapping below: newRow -> sourceRow (i.e. after the edit, row $newRow holds
# the D..T values that used to live in row $sourceRow before the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowMap = @{
    2  = 4
    3  = 8
    4  = 3
    5  = 11
    6  = 16
    7  = 5
    8  = 17
    9  = 7
    10 = 2
    11 = 12
    12 = 13
    13 = 6
    14 = 10
    15 = 19
    16 = 18
    17 = 15
    18 = 14
    19 = 9
}

# Columns that actually carry data that differs row to row (D, and M..T).
$cols = @(4, 13, 14, 15, 16, 17, 18, 19, 20)

# Snapshot the "before" values for every affected column of every row first,
# since we'll be overwriting rows in place and some rows are sources for
# others.
$snapshot = @{}
foreach ($srcRow in 2..19) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Cells.Item($srcRow, $col).Value2
    }
    $snapshot[$srcRow] = $rowVals
}

foreach ($newRow in $rowMap.Keys) {
    $srcRow = $rowMap[$newRow]
    $rowVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Cells.Item($newRow, $col).Value2 = $rowVals[$col]
    }
}

Write-Output "done"
